# mip_winners.xlsx edit:
#  - rename header "year_x" (A1) -> "season_ending_year_x"
#  - rename header "year_y" (O1) -> "season_ending_year_y"
#  - add a new trailing column AY "calendar_year" = numeric value of year_x (column A)
#  - backfill column Q "birth_year" = year_x (column A) - age_x (column D)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header renames ---
$ws.Range("A1").Value = "season_ending_year_x"
$ws.Range("O1").Value = "season_ending_year_y"

# --- new header for the appended column (match the existing header formatting) ---
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("AY1").Value = "calendar_year"

# --- backfill data rows 2..40 ---
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $yearText = $ws.Cells.Item($r, 1).Text   # column A: year_x
    $ageText  = $ws.Cells.Item($r, 4).Text   # column D: age_x

    $year = [int]$yearText
    $age  = [int]$ageText

    $birthYear    = $year - $age
    $calendarYear = $year

    $ws.Cells.Item($r, 17).Value = $birthYear     # column Q: birth_year
    $ws.Cells.Item($r, 51).Value = $calendarYear  # column AY: calendar_year
}
